# Updates cryptos list values (price/volume) per the Oct 12 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.065.52'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').Value = '2.472.34'
$ws.Range('E3').Value = '  +2.29%  '
$ws.Range('E4').Value = '  -1.01%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '577.44'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +0.35%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '146.74'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +0.90%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.42%  '
$ws.Range('D9').Value = '2.471.93'
$ws.Range('E9').Value = '  +1.01%  '
$ws.Range('E10').Value = '  +0.33%  '
$ws.Range('E11').Value = '  +1.07%  '
$ws.Range('E12').Value = '  +1.02%  '
$ws.Range('E13').Value = '  -0.23%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '28.98'
$c.Style = "Normal"
$ws.Range('E14').Value = '  +7.18%  '
$ws.Range('E15').Value = '  +0.39%  '
$ws.Range('E16').Value = '  -1.26%  '
$ws.Range('D17').Value = '63.164.15'
$ws.Range('E17').Value = '  +0.63%  '
$ws.Range('D18').Value = '2.479.70'
$ws.Range('E18').Value = '  +1.10%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '8.18'
$c.Style = "Normal"
$ws.Range('E19').Value = '  +3.54%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '11.03'
$c.Style = "Normal"
$ws.Range('E20').Value = '  +0.73%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '329.80'
$c.Style = "Normal"
$ws.Range('E21').Value = '  +0.58%  '
$ws.Range('E22').Value = '  +9.52%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('E24').Value = '  -0.09%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '66.28'
$c.Style = "Normal"
$ws.Range('E25').Value = '  +0.99%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '666.78'
$c.Style = "Normal"
$ws.Range('E26').Value = '  +7.39%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '9.63'
$c.Style = "Normal"
$ws.Range('E27').Value = '  +14.77%  '
$ws.Range('D28').Value = '0.0₃0986'
$ws.Range('E28').Value = '  +0.68%  '
$ws.Range('D29').Value = '2.591.17'
$ws.Range('E30').Value = '  +693.00%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '1.44'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +2.55%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '8.08'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -0.76%  '
$ws.Range('E33').Value = '  +0.76%  '
$ws.Range('E34').Value = '  -3.22%  '
$ws.Range('E35').Value = '  +3.50%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +0.38%  '
$ws.Range('E37').Value = '  +0.54%  '
# Rows 38-40 reorder: Monero, PolygonEcosystemToken, RenderToken
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '152.64'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('B39').Value = 'PolygonEcosystemToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.372'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -0.36%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '5.42'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +0.73%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '18.78'
$c.Style = "Normal"
$ws.Range('E41').Value = '  +0.70%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '2.73'
$c.Style = "Normal"
$ws.Range('E42').Value = '  -1.20%  '
$ws.Range('E43').Value = '  -0.39%  '
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('E45').Value = '  +6.87%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '151.87'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +5.07%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '15.14'
$c.Style = "Normal"
$ws.Range('E47').Value = '  +24.95%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '3.59'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +0.22%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '20.68'
$c.Style = "Normal"
$ws.Range('E49').Value = '  +2.28%  '
$ws.Range('E50').Value = '  +0.89%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '0.0513'
$c.Style = "Normal"
$ws.Range('E51').Value = '  -0.64%  '

Write-Output "cryptos list updated"
